# Remove the existing slide comment (the "date and time should be in the
# same line?" review comment on slide 1), as part of the Issue686 cleanup
# tweaks for the Features page screenshots.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Delete comments from the end backwards to be safe with index shifting.
for ($i = $s.Comments.Count; $i -ge 1; $i--) {
    $s.Comments.Item($i).Delete()
}
